$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.207.01'
$ws.Range('E2').Value = '  +1.20%  '
$ws.Range('D3').Value = '2.694.53'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''610.52'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '''159.88'
$ws.Range('E6').Value = '  +1.45%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '''0.591'
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('E9').Value = '  +8.37%  '
$ws.Range('D10').Value = '''6.00'
$ws.Range('E10').Value = '  +2.96%  '
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('E13').Value = '  +21.56%  '
$ws.Range('D14').Value = '''30.29'
$ws.Range('E14').Value = '  +3.61%  '
$ws.Range('D15').Value = '3.181.32'
$ws.Range('E15').Value = '  +1.56%  '
$ws.Range('D16').Value = '66.054.55'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('D17').Value = '2.686.69'
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('D18').Value = '''12.77'
$ws.Range('E18').Value = '  +1.13%  '
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('D20').Value = '''363.23'
$ws.Range('E20').Value = '  +2.32%  '
$ws.Range('D21').Value = '''7.55'
$ws.Range('E21').Value = '  +3.39%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '''70.14'
$ws.Range('E23').Value = '  +2.76%  '
$ws.Range('D24').Value = '''9.75'
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('E25').Value = '  +16.14%  '
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('D27').Value = '''0.174'
$ws.Range('E27').Value = '  +5.09%  '
$ws.Range('E28').Value = '  -0.60%  '
$ws.Range('D29').Value = '''8.26'
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('D30').Value = '''2.21'
$ws.Range('E30').Value = '  +6.14%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '''538.67'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('E33').Value = '  -1.67%  '
$ws.Range('D34').Value = '''6.62'
$ws.Range('E34').Value = '  +2.05%  '
$ws.Range('E35').Value = '  -6.16%  '
$ws.Range('D36').Value = '''0.434'
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('D37').Value = '''20.86'
$ws.Range('E37').Value = '  +3.01%  '
$ws.Range('D38').Value = '''163.27'
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('D39').Value = '''2.01'
$ws.Range('E39').Value = '  -2.42%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').Value = '''170.77'
$ws.Range('E41').Value = '  +1.06%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = '''42.70'
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('E44').Value = '  +2.61%  '
$ws.Range('E45').Value = '  +4.49%  '
$ws.Range('D46').Value = '''0.0621'
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('D47').Value = '''23.32'
$ws.Range('E47').Value = '  -0.70%  '
$ws.Range('E48').Value = '  +2.24%  '
$ws.Range('D49').Value = '''0.0267'
$ws.Range('E49').Value = '  +5.92%  '
$ws.Range('D50').Value = '''20.50'
$ws.Range('E50').Value = '  +4.75%  '
$ws.Range('E51').Value = '  +0.52%  '
